# Extra wait for loader added in user, ba, customer and withholding page
#
# This updates the QA test-data fixture row (row 2) on the
# "Memo_Verification_details" and "Memo_invoice_Details" worksheets with a
# fresh set of sample invoice values (new invoice number / date / amounts),
# matching the refreshed test fixture used by the loader-wait tests.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Memo_Verification_details  (row 2, columns B..O)
# ---------------------------------------------------------------------------
$wsVerif = $wb.Worksheets.Item("Memo_Verification_details")

# Several of these columns hold numeric-looking values that must stay text
# (shared-string) cells, exactly like the rest of the sheet, so force text
# number formatting before writing them.
$wsVerif.Range("B2:M2").NumberFormat = "@"
$wsVerif.Range("O2").NumberFormat = "@"

$wsVerif.Range("B2").Value = "TESTINV45388"
$wsVerif.Range("C2").Value = "4500000888"
$wsVerif.Range("D2").Value = "10"
$wsVerif.Range("E2").Value = "1.00"
$wsVerif.Range("F2").Value = "KC"
$wsVerif.Range("G2").Value = "996713"
$wsVerif.Range("H2").Value = "1022"
$wsVerif.Range("I2").Value = "A&E"
$wsVerif.Range("J2").Value = "996713"
$wsVerif.Range("K2").Value = "KC"
$wsVerif.Range("L2").Value = "00-00-0%"
$wsVerif.Range("M2").Value = "test"
# N2 (HSN_Code) and O2 (Tax_Code) keep their original values (996713 / KG)
$wsVerif.Range("O2").Value = "KG"

# ---------------------------------------------------------------------------
# Sheet: Memo_invoice_Details  (row 2, columns B..P)
# ---------------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item("Memo_invoice_Details")

$wsInv.Range("B2:C2").NumberFormat = "@"
$wsInv.Range("E2:P2").NumberFormat = "@"

$wsInv.Range("B2").Value = "TESTINV45388"
$wsInv.Range("C2").Value = "2024-03-15"
# D2 (Service_Name) keeps its original value
$wsInv.Range("E2").Value = "1"
$wsInv.Range("F2").Value = "0"
$wsInv.Range("G2").Value = "0"
$wsInv.Range("H2").Value = "0"
$wsInv.Range("I2").Value = "0.18"
$wsInv.Range("J2").Value = "1"
$wsInv.Range("K2").Value = "1"
$wsInv.Range("L2").Value = "996713"
$wsInv.Range("M2").Value = "test_customer_1"
$wsInv.Range("N2").Value = "test_comment_1"
$wsInv.Range("O2").Value = "1.18"
$wsInv.Range("P2").Value = "1022"
